$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Uppercase a couple of header labels
$ws.Range("I1").Value = "TOTAL PHYSICAL TARGET"
$ws.Range("L1").Value = "BATCH"

# 2. Insert 5 new columns before the old "Status as of ..." column (AA),
#    pushing it (and its data validation) out to AF. Insert() on whole
#    columns preserves the existing AA column's contents/formatting by
#    shifting it right, and creates new blank AA:AE columns that inherit
#    the header row's style from the column that used to be there.
$ws.Range("AA1:AE1").EntireColumn.Insert()

# 3. Populate the five new header cells (AA1:AE1)
$ws.Range("AA1").Value = "No. of Sites Reverted"
$ws.Range("AB1").Value = "No. of Sites Not yet started"
$ws.Range("AC1").Value = "No. of Sites Under Procurement"
$ws.Range("AD1").Value = "No. of Sites On Going"
$ws.Range("AE1").Value = "No. of Sites Completed"

# 4. Clear out the placeholder "-" entries that used to live in the
#    "Total Physical Target" (I) and "Batch" (L) columns for rows 2-48.
$ws.Range("I2:I48").ClearContents()
$ws.Range("L2:L48").ClearContents()
